$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.745.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.64%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.759.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "

# Row 6
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4444"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3736"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.84%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07734"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.128"
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.77%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.197"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.69%  "

# Row 15
$ws.Range("E15").Value = "  -1.86%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.759.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.72%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +13.56%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001081"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.96%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06239"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.77%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.199"
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5326"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.772.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.55%  "

# Row 25
$ws.Range("E25").Value = "  -1.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.323"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.40%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.364"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.959.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.75%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.82%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.217"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.786"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.20%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09286"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.62%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.657"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.46%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.71"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.95%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02345"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.57%  "

# Row 38
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2192"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.64%  "

# Row 39
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6513"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.94%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06156"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.099"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.82%  "

# Row 42
$ws.Range("E42").Value = "  -0.46%  "

# Row 43
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.431"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.36%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.031"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.79%  "

# Row 45
$ws.Range("E45").Value = "  +0.03%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.90%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6038"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.15%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.764"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.64%  "

# Row 50
$ws.Range("E50").Value = "  -1.29%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.148"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.00%  "
